# notebook 2 update FIXED BUG
# Rescale the "Calculated Ribeye Area" (column D) and
# "Calculated Fat Thickness" (column E) values by the corrected
# model factor of 281/280 for all data rows (2-41).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$factor = 281.0 / 280.0

$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)   # Column D
    $eCell = $ws.Cells.Item($r, 5)   # Column E

    $dCell.Value2 = $dCell.Value2 * $factor
    $eCell.Value2 = $eCell.Value2 * $factor
}
